# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" list (col E, rows 16-38) was loaded most-recent-first
# (2412 down to 2302). Re-sort it to read oldest -> newest (2302 up to
# 2412), carrying each period's "Valor Mora" (col F) along with it so the
# figures stay attached to the correct period.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 16
$lastRow = 38

# Snapshot the current period / mora-value pairs.
$periods = @()
$moraValues = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $periods += , $ws.Cells.Item($r, 5).Value2
    $moraValues += , $ws.Cells.Item($r, 6).Value2
}

# Write the pairs back in reverse order (ascending chronological order).
$n = $periods.Length
for ($i = 0; $i -lt $n; $i++) {
    $r = $firstRow + $i
    $ws.Cells.Item($r, 5).Value = $periods[$n - 1 - $i]
    $ws.Cells.Item($r, 6).Value = $moraValues[$n - 1 - $i]
}
